$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column D (ag_non_ag_group) before existing lu_group column
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range('A1').Value = 'code'
$ws.Range('B1').Value = 'desc'
$ws.Range('C1').Value = 'ag_group'
$ws.Range('D1').Value = 'ag_non_ag_group'
$ws.Range('E1').Value = 'lu_group'

# Main land-use rows (codes -2..27, 100..107)
$ws.Cells.Item(2,1).Value = -2
$ws.Cells.Item(2,2).Value = 'Other land-use'
$ws.Cells.Item(2,5).Value = 'Other land-use'
$ws.Cells.Item(3,1).Value = -1
$ws.Cells.Item(3,2).Value = 'Other land-use'
$ws.Cells.Item(3,5).Value = 'Other land-use'
$ws.Cells.Item(4,1).Value = 0
$ws.Cells.Item(4,2).Value = 'Apples'
$ws.Cells.Item(4,3).Value = 'Crops'
$ws.Cells.Item(4,4).Value = 'Crops'
$ws.Cells.Item(4,5).Value = 'Agricultural land-use'
$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = 'Beef - modified land'
$ws.Cells.Item(5,3).Value = 'Livestock'
$ws.Cells.Item(5,4).Value = 'Livestock'
$ws.Cells.Item(5,5).Value = 'Agricultural land-use'
$ws.Cells.Item(6,1).Value = 2
$ws.Cells.Item(6,2).Value = 'Beef - natural land'
$ws.Cells.Item(6,3).Value = 'Livestock'
$ws.Cells.Item(6,4).Value = 'Livestock'
$ws.Cells.Item(6,5).Value = 'Agricultural land-use'
$ws.Cells.Item(7,1).Value = 3
$ws.Cells.Item(7,2).Value = 'Citrus'
$ws.Cells.Item(7,3).Value = 'Crops'
$ws.Cells.Item(7,4).Value = 'Crops'
$ws.Cells.Item(7,5).Value = 'Agricultural land-use'
$ws.Cells.Item(8,1).Value = 4
$ws.Cells.Item(8,2).Value = 'Cotton'
$ws.Cells.Item(8,3).Value = 'Crops'
$ws.Cells.Item(8,4).Value = 'Crops'
$ws.Cells.Item(8,5).Value = 'Agricultural land-use'
$ws.Cells.Item(9,1).Value = 5
$ws.Cells.Item(9,2).Value = 'Dairy - modified land'
$ws.Cells.Item(9,3).Value = 'Livestock'
$ws.Cells.Item(9,4).Value = 'Livestock'
$ws.Cells.Item(9,5).Value = 'Agricultural land-use'
$ws.Cells.Item(10,1).Value = 6
$ws.Cells.Item(10,2).Value = 'Dairy - natural land'
$ws.Cells.Item(10,3).Value = 'Livestock'
$ws.Cells.Item(10,4).Value = 'Livestock'
$ws.Cells.Item(10,5).Value = 'Agricultural land-use'
$ws.Cells.Item(11,1).Value = 7
$ws.Cells.Item(11,2).Value = 'Grapes'
$ws.Cells.Item(11,3).Value = 'Crops'
$ws.Cells.Item(11,4).Value = 'Crops'
$ws.Cells.Item(11,5).Value = 'Agricultural land-use'
$ws.Cells.Item(12,1).Value = 8
$ws.Cells.Item(12,2).Value = 'Hay'
$ws.Cells.Item(12,3).Value = 'Crops'
$ws.Cells.Item(12,4).Value = 'Crops'
$ws.Cells.Item(12,5).Value = 'Agricultural land-use'
$ws.Cells.Item(13,1).Value = 9
$ws.Cells.Item(13,2).Value = 'Nuts'
$ws.Cells.Item(13,3).Value = 'Crops'
$ws.Cells.Item(13,4).Value = 'Crops'
$ws.Cells.Item(13,5).Value = 'Agricultural land-use'
$ws.Cells.Item(14,1).Value = 10
$ws.Cells.Item(14,2).Value = 'Other non-cereal crops'
$ws.Cells.Item(14,3).Value = 'Crops'
$ws.Cells.Item(14,4).Value = 'Crops'
$ws.Cells.Item(14,5).Value = 'Agricultural land-use'
$ws.Cells.Item(15,1).Value = 11
$ws.Cells.Item(15,2).Value = 'Pears'
$ws.Cells.Item(15,3).Value = 'Crops'
$ws.Cells.Item(15,4).Value = 'Crops'
$ws.Cells.Item(15,5).Value = 'Agricultural land-use'
$ws.Cells.Item(16,1).Value = 12
$ws.Cells.Item(16,2).Value = 'Plantation fruit'
$ws.Cells.Item(16,3).Value = 'Crops'
$ws.Cells.Item(16,4).Value = 'Crops'
$ws.Cells.Item(16,5).Value = 'Agricultural land-use'
$ws.Cells.Item(17,1).Value = 13
$ws.Cells.Item(17,2).Value = 'Rice'
$ws.Cells.Item(17,3).Value = 'Crops'
$ws.Cells.Item(17,4).Value = 'Crops'
$ws.Cells.Item(17,5).Value = 'Agricultural land-use'
$ws.Cells.Item(18,1).Value = 14
$ws.Cells.Item(18,2).Value = 'Sheep - modified land'
$ws.Cells.Item(18,3).Value = 'Livestock'
$ws.Cells.Item(18,4).Value = 'Livestock'
$ws.Cells.Item(18,5).Value = 'Agricultural land-use'
$ws.Cells.Item(19,1).Value = 15
$ws.Cells.Item(19,2).Value = 'Sheep - natural land'
$ws.Cells.Item(19,3).Value = 'Livestock'
$ws.Cells.Item(19,4).Value = 'Livestock'
$ws.Cells.Item(19,5).Value = 'Agricultural land-use'
$ws.Cells.Item(20,1).Value = 16
$ws.Cells.Item(20,2).Value = 'Stone fruit'
$ws.Cells.Item(20,3).Value = 'Crops'
$ws.Cells.Item(20,4).Value = 'Crops'
$ws.Cells.Item(20,5).Value = 'Agricultural land-use'
$ws.Cells.Item(21,1).Value = 17
$ws.Cells.Item(21,2).Value = 'Sugar'
$ws.Cells.Item(21,3).Value = 'Crops'
$ws.Cells.Item(21,4).Value = 'Crops'
$ws.Cells.Item(21,5).Value = 'Agricultural land-use'
$ws.Cells.Item(22,1).Value = 18
$ws.Cells.Item(22,2).Value = 'Summer cereals'
$ws.Cells.Item(22,3).Value = 'Crops'
$ws.Cells.Item(22,4).Value = 'Crops'
$ws.Cells.Item(22,5).Value = 'Agricultural land-use'
$ws.Cells.Item(23,1).Value = 19
$ws.Cells.Item(23,2).Value = 'Summer legumes'
$ws.Cells.Item(23,3).Value = 'Crops'
$ws.Cells.Item(23,4).Value = 'Crops'
$ws.Cells.Item(23,5).Value = 'Agricultural land-use'
$ws.Cells.Item(24,1).Value = 20
$ws.Cells.Item(24,2).Value = 'Summer oilseeds'
$ws.Cells.Item(24,3).Value = 'Crops'
$ws.Cells.Item(24,4).Value = 'Crops'
$ws.Cells.Item(24,5).Value = 'Agricultural land-use'
$ws.Cells.Item(25,1).Value = 21
$ws.Cells.Item(25,2).Value = 'Tropical stone fruit'
$ws.Cells.Item(25,3).Value = 'Crops'
$ws.Cells.Item(25,4).Value = 'Crops'
$ws.Cells.Item(25,5).Value = 'Agricultural land-use'
$ws.Cells.Item(26,1).Value = 22
$ws.Cells.Item(26,2).Value = 'Unallocated - modified land'
$ws.Cells.Item(26,3).Value = 'Unallocated - modified land'
$ws.Cells.Item(26,4).Value = 'Unallocated - modified land'
$ws.Cells.Item(26,5).Value = 'Agricultural land-use'
$ws.Cells.Item(27,1).Value = 23
$ws.Cells.Item(27,2).Value = 'Unallocated - natural land'
$ws.Cells.Item(27,3).Value = 'Unallocated - natural land'
$ws.Cells.Item(27,4).Value = 'Unallocated - natural land'
$ws.Cells.Item(27,5).Value = 'Agricultural land-use'
$ws.Cells.Item(28,1).Value = 24
$ws.Cells.Item(28,2).Value = 'Vegetables'
$ws.Cells.Item(28,3).Value = 'Crops'
$ws.Cells.Item(28,4).Value = 'Crops'
$ws.Cells.Item(28,5).Value = 'Agricultural land-use'
$ws.Cells.Item(29,1).Value = 25
$ws.Cells.Item(29,2).Value = 'Winter cereals'
$ws.Cells.Item(29,3).Value = 'Crops'
$ws.Cells.Item(29,4).Value = 'Crops'
$ws.Cells.Item(29,5).Value = 'Agricultural land-use'
$ws.Cells.Item(30,1).Value = 26
$ws.Cells.Item(30,2).Value = 'Winter legumes'
$ws.Cells.Item(30,3).Value = 'Crops'
$ws.Cells.Item(30,4).Value = 'Crops'
$ws.Cells.Item(30,5).Value = 'Agricultural land-use'
$ws.Cells.Item(31,1).Value = 27
$ws.Cells.Item(31,2).Value = 'Winter oilseeds'
$ws.Cells.Item(31,3).Value = 'Crops'
$ws.Cells.Item(31,4).Value = 'Crops'
$ws.Cells.Item(31,5).Value = 'Agricultural land-use'
$ws.Cells.Item(32,1).Value = 100
$ws.Cells.Item(32,2).Value = 'Environmental Plantings'
$ws.Cells.Item(32,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(32,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(33,1).Value = 101
$ws.Cells.Item(33,2).Value = 'Riparian Plantings'
$ws.Cells.Item(33,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(33,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(34,1).Value = 102
$ws.Cells.Item(34,2).Value = 'Sheep Agroforestry'
$ws.Cells.Item(34,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(34,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(35,1).Value = 103
$ws.Cells.Item(35,2).Value = 'Beef Agroforestry'
$ws.Cells.Item(35,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(35,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(36,1).Value = 104
$ws.Cells.Item(36,2).Value = 'Carbon Plantings (Block)'
$ws.Cells.Item(36,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(36,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(37,1).Value = 105
$ws.Cells.Item(37,2).Value = 'Sheep Carbon Plantings (Belt)'
$ws.Cells.Item(37,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(37,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(38,1).Value = 106
$ws.Cells.Item(38,2).Value = 'Beef Carbon Plantings (Belt)'
$ws.Cells.Item(38,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(38,5).Value = 'Non-agricultural land-use'
$ws.Cells.Item(39,1).Value = 107
$ws.Cells.Item(39,2).Value = 'BECCS'
$ws.Cells.Item(39,4).Value = 'Non-agricultural land-use'
$ws.Cells.Item(39,5).Value = 'Non-agricultural land-use'

# Agricultural management rows
$ws.Cells.Item(40,1).Value = 1
$ws.Cells.Item(40,2).Value = 'Asparagopsis taxiformis'
$ws.Cells.Item(40,5).Value = 'Agricultural management'
$ws.Cells.Item(41,1).Value = 2
$ws.Cells.Item(41,2).Value = 'Precision Agriculture'
$ws.Cells.Item(41,5).Value = 'Agricultural management'
$ws.Cells.Item(42,1).Value = 3
$ws.Cells.Item(42,2).Value = 'Ecological Grazing'
$ws.Cells.Item(42,5).Value = 'Agricultural management'
$ws.Cells.Item(43,1).Value = 4
$ws.Cells.Item(43,2).Value = 'Savanna Burning'
$ws.Cells.Item(43,5).Value = 'Agricultural management'
$ws.Cells.Item(44,1).Value = 5
$ws.Cells.Item(44,2).Value = 'AgTech EI'
$ws.Cells.Item(44,5).Value = 'Agricultural management'
$ws.Cells.Item(45,2).Value = 'Deforestation'
$ws.Cells.Item(45,5).Value = 'Transition'

# Column B width
$ws.Columns.Item(2).ColumnWidth = 27

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# View / selection
$ws.Range("B27").Select()
